$wb = $excel.ActiveWorkbook

# Scheduled market-data refresh: update currentAveragePrice*, LevePrice*,
# and LeveProfit* columns (H-N) across the profession Leve-profit tables.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4200
$ws.Range("I64").Value = 3800
$ws.Range("J64").Value = 4314.2856
$ws.Range("K64").Value = 3800
$ws.Range("L64").Value = 4314.2856
$ws.Range("M64").Value = -3552
$ws.Range("N64").Value = -4810.2856

$ws.Range("H67").Value = 4200
$ws.Range("I67").Value = 3800
$ws.Range("J67").Value = 4314.2856
$ws.Range("K67").Value = 3800
$ws.Range("L67").Value = 4314.2856
$ws.Range("M67").Value = -2942
$ws.Range("N67").Value = -6030.2856

$ws.Range("H137").Value = 1853215.8
$ws.Range("I137").Value = 2501045
$ws.Range("K137").Value = 7503135
$ws.Range("M137").Value = -7500585

$ws.Range("H138").Value = 2224360.5
$ws.Range("I138").Value = 1463.4
$ws.Range("J138").Value = 4169395.5
$ws.Range("K138").Value = 4390.200000000001
$ws.Range("L138").Value = 12508186.5
$ws.Range("M138").Value = 749.7999999999993
$ws.Range("N138").Value = -12518466.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2811.8262
$ws.Range("I2").Value = 1941
$ws.Range("J2").Value = 3481.6924
$ws.Range("K2").Value = 1941
$ws.Range("L2").Value = 3481.6924
$ws.Range("M2").Value = -1828
$ws.Range("N2").Value = -3707.6924

$ws.Range("H26").Value = 7713.857
$ws.Range("I26").Value = 6499.25
$ws.Range("J26").Value = 9333.333000000001
$ws.Range("K26").Value = 6499.25
$ws.Range("L26").Value = 9333.333000000001
$ws.Range("M26").Value = -6169.25
$ws.Range("N26").Value = -9993.333000000001

$ws.Range("H32").Value = 364.97
$ws.Range("I32").Value = 370.957
$ws.Range("J32").Value = 285.42856
$ws.Range("K32").Value = 370.957
$ws.Range("L32").Value = 285.42856
$ws.Range("M32").Value = -83.95699999999999
$ws.Range("N32").Value = -859.4285600000001

$ws.Range("H74").Value = 9334681
$ws.Range("J74").Value = 92628.55
$ws.Range("L74").Value = 92628.55
$ws.Range("N74").Value = -94376.55

$ws.Range("H77").Value = 9334681
$ws.Range("J77").Value = 92628.55
$ws.Range("L77").Value = 463142.75
$ws.Range("N77").Value = -471878.75

$ws.Range("H116").Value = 2811.8262
$ws.Range("I116").Value = 1941
$ws.Range("J116").Value = 3481.6924
$ws.Range("K116").Value = 1941
$ws.Range("L116").Value = 3481.6924
$ws.Range("M116").Value = 353
$ws.Range("N116").Value = -8069.6924

$ws.Range("H132").Value = 36189.934
$ws.Range("I132").Value = 21578.104
$ws.Range("J132").Value = 99950.63
$ws.Range("K132").Value = 64734.312
$ws.Range("L132").Value = 299851.89
$ws.Range("M132").Value = -62204.312
$ws.Range("N132").Value = -304911.89

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2811.8262
$ws.Range("I3").Value = 1941
$ws.Range("J3").Value = 3481.6924
$ws.Range("K3").Value = 1941
$ws.Range("L3").Value = 3481.6924
$ws.Range("M3").Value = -1827
$ws.Range("N3").Value = -3709.6924

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H134").Value = 1547.7819
$ws.Range("I134").Value = 1002.2683
$ws.Range("J134").Value = 3145.3572
$ws.Range("K134").Value = 3006.8049
$ws.Range("L134").Value = 9436.071599999999
$ws.Range("M134").Value = -471.8049000000001
$ws.Range("N134").Value = -14506.0716

$ws.Range("H141").Value = 49244
$ws.Range("I141").Value = 40709
$ws.Range("J141").Value = 52658
$ws.Range("K141").Value = 40709
$ws.Range("L141").Value = 52658
$ws.Range("M141").Value = -35529
$ws.Range("N141").Value = -63018

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2944.575
$ws.Range("I31").Value = 1277.6786
$ws.Range("J31").Value = 6834
$ws.Range("K31").Value = 1277.6786
$ws.Range("L31").Value = 6834
$ws.Range("M31").Value = -982.6786
$ws.Range("N31").Value = -7424

$ws.Range("H34").Value = 2944.575
$ws.Range("I34").Value = 1277.6786
$ws.Range("J34").Value = 6834
$ws.Range("K34").Value = 1277.6786
$ws.Range("L34").Value = 6834
$ws.Range("M34").Value = -1075.6786
$ws.Range("N34").Value = -7238

$ws.Range("H58").Value = 21740698
$ws.Range("I58").Value = 25001444
$ws.Range("J58").Value = 2389.1667
$ws.Range("K58").Value = 25001444
$ws.Range("L58").Value = 2389.1667
$ws.Range("M58").Value = -25001241
$ws.Range("N58").Value = -2795.1667

$ws.Range("H136").Value = 21740698
$ws.Range("I136").Value = 25001444
$ws.Range("J136").Value = 2389.1667
$ws.Range("K136").Value = 75004332
$ws.Range("L136").Value = 7167.500100000001
$ws.Range("M136").Value = -75001782
$ws.Range("N136").Value = -12267.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 412.65625
$ws.Range("I5").Value = 337.22223
$ws.Range("J5").Value = 820
$ws.Range("K5").Value = 1011.66669
$ws.Range("L5").Value = 2460
$ws.Range("M5").Value = -899.66669
$ws.Range("N5").Value = -2684

$ws.Range("H135").Value = 412.65625
$ws.Range("I135").Value = 337.22223
$ws.Range("J135").Value = 820
$ws.Range("K135").Value = 3035.00007
$ws.Range("L135").Value = 7380
$ws.Range("M135").Value = -500.0000700000001
$ws.Range("N135").Value = -12450

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 85975.586
$ws.Range("I132").Value = 2929.7144
$ws.Range("K132").Value = 8789.143199999999
$ws.Range("M132").Value = -6259.143199999999

$ws.Range("H136").Value = 84900.36
$ws.Range("I136").Value = 63407.125
$ws.Range("J136").Value = 123110.555
$ws.Range("K136").Value = 190221.375
$ws.Range("L136").Value = 369331.665
$ws.Range("M136").Value = -187671.375
$ws.Range("N136").Value = -374431.665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 533.0909
$ws.Range("I107").Value = 488.33334
$ws.Range("J107").Value = 586.8
$ws.Range("K107").Value = 1465.00002
$ws.Range("L107").Value = 1760.4
$ws.Range("M107").Value = 454.9999800000001
$ws.Range("N107").Value = -5600.4

$ws.Range("H132").Value = 28854.492
$ws.Range("I132").Value = 17299.533
$ws.Range("J132").Value = 91881.55
$ws.Range("K132").Value = 51898.599
$ws.Range("L132").Value = 275644.65
$ws.Range("M132").Value = -49368.599
$ws.Range("N132").Value = -280704.65
